# B1--and-B2-PowerPoint.pptx edit
#
# 1) The table on slide 5 gets a different built-in table style applied
#    (PowerPoint Table Design gallery -> a new style GUID is written to
#    <a:tableStyleId>).
# 2) The deck's theme is swapped from the "Integral / Red Violet" theme
#    to the stock "Office Theme" colours (Design tab -> Themes gallery).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Table style on slide 5
# ---------------------------------------------------------------------
$tableSlide = $p.Slides.Item(5)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{598BE393-C4D6-4CF8-AC21-DD8676420148}")
    }
}

# ---------------------------------------------------------------------
# 2) Re-theme the deck: Integral/Red Violet -> Office Theme
# ---------------------------------------------------------------------
function Set-ThemeRGB($scheme, $index, $r, $g, $b) {
    $scheme.Item($index).RGB = $r + ($g * 256) + ($b * 65536)
}

$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

Set-ThemeRGB $colors 1  0x00 0x00 0x00   # dk1
Set-ThemeRGB $colors 2  0xFF 0xFF 0xFF   # lt1
Set-ThemeRGB $colors 3  0x44 0x54 0x6A   # dk2
Set-ThemeRGB $colors 4  0xE7 0xE6 0xE6   # lt2
Set-ThemeRGB $colors 5  0x5B 0x9B 0xD5   # accent1
Set-ThemeRGB $colors 6  0xED 0x7D 0x31   # accent2
Set-ThemeRGB $colors 7  0xA5 0xA5 0xA5   # accent3
Set-ThemeRGB $colors 8  0xFF 0xC0 0x00   # accent4
Set-ThemeRGB $colors 9  0x44 0x72 0xC4   # accent5
Set-ThemeRGB $colors 10 0x70 0xAD 0x47   # accent6
Set-ThemeRGB $colors 11 0x05 0x63 0xC1   # hlink
Set-ThemeRGB $colors 12 0x95 0x4F 0x72   # folHlink

# Theme fonts (Office Theme keeps Arial for both major/minor Latin, same
# as the previous theme, but set explicitly for completeness).
$theme.ThemeFontScheme.MajorFont.Latin = "Arial"
$theme.ThemeFontScheme.MinorFont.Latin = "Arial"
